$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value of 45180 (2023-09-11)
# for every data row (2 through 302). The update bumps it by one day to
# 45181 (2023-09-12) across the whole column range.
$lastRow = 302

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
